# Insert a new weekly record for "Macroferia Regional de Talca - Apio" at row 261.
# This pushes the existing rows 261-379 down to 262-380 (dimension grows to A1:R380)
# and populates the newly opened row 261 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 261, shifting rows 261:379 down to 262:380.
$ws.Rows(261).Insert()

# Fill the new row 261 with the new record's data.
$ws.Cells.Item(261, 1).Value = 5
$ws.Cells.Item(261, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(261, 3).Value = "Maule"
$ws.Cells.Item(261, 4).Value = 45205
$ws.Cells.Item(261, 5).Value = 7
$ws.Cells.Item(261, 6).Value = 100112017
$ws.Cells.Item(261, 7).Value = "Apio"
$ws.Cells.Item(261, 8).Value = "Americana (o)"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 700
$ws.Cells.Item(261, 11).Value = 6000
$ws.Cells.Item(261, 12).Value = 6000
$ws.Cells.Item(261, 13).Value = 6000
$ws.Cells.Item(261, 14).Value = "$/docena de matas"
$ws.Cells.Item(261, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(261, 16).Value = 1000
$ws.Cells.Item(261, 17).Value = 6
$ws.Cells.Item(261, 18).Value = "Hortaliza"
